$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9197401404380798
$ws.Range("B1").Value = 1.658546209335327
$ws.Range("C1").Value = 4.412744998931885
$ws.Range("D1").Value = 2.348679065704346
$ws.Range("E1").Value = 0.8666198253631592
